$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 10:42"

$ws.Cells.Item(9, 1).Value = "Iran"
$ws.Cells.Item(9, 2).Value = 32332
$ws.Cells.Item(9, 3).Value = 2926
$ws.Cells.Item(9, 4).Value = 11133
$ws.Cells.Item(9, 5).Value = 18821
$ws.Cells.Item(9, 6).Value = 2746
$ws.Cells.Item(9, 7).Value = 144
$ws.Cells.Item(9, 8).Value = 2378

$ws.Cells.Item(15, 1).Value = "Austria"
$ws.Cells.Item(15, 2).Value = 7129
$ws.Cells.Item(15, 3).Value = 220
$ws.Cells.Item(15, 4).Value = 225
$ws.Cells.Item(15, 5).Value = 6846
$ws.Cells.Item(15, 6).Value = 96
$ws.Cells.Item(15, 7).Value = 9
$ws.Cells.Item(15, 8).Value = 58

$ws.Cells.Item(20, 1).Value = "Noruega"
$ws.Cells.Item(20, 2).Value = 3380
$ws.Cells.Item(20, 3).Value = 8
$ws.Cells.Item(20, 4).Value = 6
$ws.Cells.Item(20, 5).Value = 3359
$ws.Cells.Item(20, 6).Value = 70
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(20, 8).Value = 15

$ws.Cells.Item(25, 1).Value = "Malasia"
$ws.Cells.Item(25, 2).Value = 2161
$ws.Cells.Item(25, 3).Value = 130
$ws.Cells.Item(25, 4).Value = 259
$ws.Cells.Item(25, 5).Value = 1876
$ws.Cells.Item(25, 6).Value = 54
$ws.Cells.Item(25, 7).Value = 3
$ws.Cells.Item(25, 8).Value = 26

$ws.Cells.Item(26, 1).Value = "Chequia"
$ws.Cells.Item(26, 2).Value = 2062
$ws.Cells.Item(26, 3).Value = 137
$ws.Cells.Item(26, 4).Value = 10
$ws.Cells.Item(26, 5).Value = 2043
$ws.Cells.Item(26, 6).Value = 34
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 9

$ws.Cells.Item(27, 1).Value = "Dinamarca"
$ws.Cells.Item(27, 2).Value = 2010
$ws.Cells.Item(27, 3).Value = 133
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 1968
$ws.Cells.Item(27, 6).Value = 94
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 41

$ws.Cells.Item(40, 1).Value = "Finlandia"
$ws.Cells.Item(40, 2).Value = 1004
$ws.Cells.Item(40, 3).Value = 46
$ws.Cells.Item(40, 4).Value = 10
$ws.Cells.Item(40, 5).Value = 989
$ws.Cells.Item(40, 6).Value = 24
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 5

$ws.Cells.Item(45, 1).Value = "India"
$ws.Cells.Item(45, 2).Value = 761
$ws.Cells.Item(45, 3).Value = 34
$ws.Cells.Item(45, 4).Value = 71
$ws.Cells.Item(45, 5).Value = 670
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 20

$ws.Cells.Item(52, 1).Value = "Estonia"
$ws.Cells.Item(52, 2).Value = 575
$ws.Cells.Item(52, 3).Value = 37
$ws.Cells.Item(52, 4).Value = 8
$ws.Cells.Item(52, 5).Value = 566
$ws.Cells.Item(52, 6).Value = 6
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 1

$ws.Cells.Item(53, 1).Value = "Eslovenia"
$ws.Cells.Item(53, 2).Value = 562
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 10
$ws.Cells.Item(53, 5).Value = 546
$ws.Cells.Item(53, 6).Value = 14
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 6

$ws.Cells.Item(54, 1).Value = "Croacia"
$ws.Cells.Item(54, 2).Value = 551
$ws.Cells.Item(54, 3).Value = 56
$ws.Cells.Item(54, 4).Value = 37
$ws.Cells.Item(54, 5).Value = 511
$ws.Cells.Item(54, 6).Value = 14
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 3

$ws.Cells.Item(55, 1).Value = "Catar"
$ws.Cells.Item(55, 2).Value = 549
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 43
$ws.Cells.Item(55, 5).Value = 506
$ws.Cells.Item(55, 6).Value = 6
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 0

$ws.Cells.Item(76, 1).Value = "Tunez"
$ws.Cells.Item(76, 2).Value = 227
$ws.Cells.Item(76, 3).Value = 30
$ws.Cells.Item(76, 4).Value = 2
$ws.Cells.Item(76, 5).Value = 219
$ws.Cells.Item(76, 6).Value = 10
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 6

$ws.Cells.Item(77, 1).Value = "Eslovaquia"
$ws.Cells.Item(77, 2).Value = 226
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 2
$ws.Cells.Item(77, 5).Value = 224
$ws.Cells.Item(77, 6).Value = 2
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0

$ws.Cells.Item(78, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(78, 2).Value = 226
$ws.Cells.Item(78, 3).Value = 35
$ws.Cells.Item(78, 4).Value = 5
$ws.Cells.Item(78, 5).Value = 218
$ws.Cells.Item(78, 6).Value = 1
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 3

$ws.Cells.Item(79, 1).Value = "Principado de Andorra"
$ws.Cells.Item(79, 2).Value = 224
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(79, 4).Value = 1
$ws.Cells.Item(79, 5).Value = 220
$ws.Cells.Item(79, 6).Value = 6
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 3

$ws.Cells.Item(80, 1).Value = "Ucrania"
$ws.Cells.Item(80, 2).Value = 218
$ws.Cells.Item(80, 3).Value = 22
$ws.Cells.Item(80, 4).Value = 4
$ws.Cells.Item(80, 5).Value = 209
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 5

$ws.Cells.Item(81, 1).Value = "Jordania"
$ws.Cells.Item(81, 2).Value = 212
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = 1
$ws.Cells.Item(81, 5).Value = 211
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 0

$ws.Cells.Item(82, 1).Value = "San Marino"
$ws.Cells.Item(82, 2).Value = 208
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = 4
$ws.Cells.Item(82, 5).Value = 183
$ws.Cells.Item(82, 6).Value = 12
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 21

$ws.Cells.Item(83, 1).Value = "Kuwait"
$ws.Cells.Item(83, 2).Value = 208
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 49
$ws.Cells.Item(83, 5).Value = 159
$ws.Cells.Item(83, 6).Value = 7
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 0

$ws.Cells.Item(84, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(84, 2).Value = 201
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = 3
$ws.Cells.Item(84, 5).Value = 195
$ws.Cells.Item(84, 6).Value = 1
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 3

$ws.Cells.Item(94, 1).Value = "Oman"
$ws.Cells.Item(94, 2).Value = 131
$ws.Cells.Item(94, 3).Value = 22
$ws.Cells.Item(94, 4).Value = 23
$ws.Cells.Item(94, 5).Value = 108
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 0

$ws.Cells.Item(95, 1).Value = "Kazajistan"
$ws.Cells.Item(95, 2).Value = 125
$ws.Cells.Item(95, 3).Value = 12
$ws.Cells.Item(95, 4).Value = 2
$ws.Cells.Item(95, 5).Value = 122
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 1

$ws.Cells.Item(96, 1).Value = "Azerbaiyan"
$ws.Cells.Item(96, 2).Value = 122
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 15
$ws.Cells.Item(96, 5).Value = 104
$ws.Cells.Item(96, 6).Value = 6
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 3

$ws.Cells.Item(97, 1).Value = "Brunei"
$ws.Cells.Item(97, 2).Value = 114
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 5
$ws.Cells.Item(97, 5).Value = 109
$ws.Cells.Item(97, 6).Value = 1
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0

$ws.Cells.Item(104, 1).Value = "Estado de Palestina"
$ws.Cells.Item(104, 2).Value = 91
$ws.Cells.Item(104, 3).Value = 5
$ws.Cells.Item(104, 4).Value = 17
$ws.Cells.Item(104, 5).Value = 73
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 1

$ws.Cells.Item(105, 1).Value = "Camerun"
$ws.Cells.Item(105, 2).Value = 88
$ws.Cells.Item(105, 3).Value = 13
$ws.Cells.Item(105, 4).Value = 2
$ws.Cells.Item(105, 5).Value = 84
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = 2

$ws.Cells.Item(143, 1).Value = "Etiopia"
$ws.Cells.Item(143, 2).Value = 16
$ws.Cells.Item(143, 3).Value = 4
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 16
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

$ws.Cells.Item(144, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(144, 2).Value = 15
$ws.Cells.Item(144, 3).Value = 1
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 15
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 0

$ws.Cells.Item(145, 1).Value = "Bermudas"
$ws.Cells.Item(145, 2).Value = 15
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = 13
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 0

$ws.Cells.Item(147, 1).Value = "El Salvador"
$ws.Cells.Item(147, 2).Value = 13
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 0
$ws.Cells.Item(147, 5).Value = 13
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 0

$ws.Cells.Item(148, 1).Value = "Maldivas"
$ws.Cells.Item(148, 2).Value = 13
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 9
$ws.Cells.Item(148, 5).Value = 4
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 0

$ws.Cells.Item(150, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(150, 2).Value = 11
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 0
$ws.Cells.Item(150, 5).Value = 11
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 0

$ws.Cells.Item(151, 1).Value = "Dominica"
$ws.Cells.Item(151, 2).Value = 11
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 11
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 0

$ws.Cells.Item(152, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(152, 2).Value = 11
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 0
$ws.Cells.Item(152, 5).Value = 11
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 0

$ws.Cells.Item(153, 1).Value = "Mongolia"
$ws.Cells.Item(153, 2).Value = 11
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 11
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 0

$ws.Cells.Item(157, 1).Value = "Haiti"
$ws.Cells.Item(157, 2).Value = 8
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 0
$ws.Cells.Item(157, 5).Value = 8
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 0

$ws.Cells.Item(158, 1).Value = "Surinam"
$ws.Cells.Item(158, 2).Value = 8
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 8
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 0
